$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 18.94467633333333
$ws.Range("H2").Value = 56.834029
$ws.Range("I2").Value = 0.03347881112463321
$ws.Range("J2").Value = 0.03347881112463321
$ws.Range("M2").Value = 15.03663066666667
$ws.Range("N2").Value = 45.109892
$ws.Range("O2").Value = 0.279146411176606
$ws.Range("P2").Value = 0.279146411176606
$ws.Range("Q2").Value = 284.8641011238743
$ws.Range("R2").Value = 2563.776910114868
$ws.Range("S2").Value = 0.009345489975900794
$ws.Range("T2").Value = 0.009345489975900793
$ws.Range("G3").Value = 18.94467633333333
$ws.Range("H3").Value = 56.834029
$ws.Range("I3").Value = 0.03347881112463321
$ws.Range("J3").Value = 0.03347881112463321
$ws.Range("O3").Value = 0.2673306493381863
$ws.Range("P3").Value = 0.2673306493381863
$ws.Range("Q3").Value = 272.8063198290765
$ws.Range("R3").Value = 2455.256878461688
$ws.Range("S3").Value = 0.008949912317018692
$ws.Range("T3").Value = 0.008949912317018692
$ws.Range("G4").Value = 18.94467633333333
$ws.Range("H4").Value = 56.834029
$ws.Range("I4").Value = 0.03347881112463321
$ws.Range("J4").Value = 0.03347881112463321
$ws.Range("M4").Value = 22.16851266666667
$ws.Range("N4").Value = 66.505538
$ws.Range("O4").Value = 0.411545703901694
$ws.Range("P4").Value = 0.411545703901694
$ws.Range("Q4").Value = 419.9752972614003
$ws.Range("R4").Value = 3779.777675352602
$ws.Range("S4").Value = 0.01377806089007904
$ws.Range("T4").Value = 0.01377806089007904
$ws.Range("G5").Value = 18.94467633333333
$ws.Range("H5").Value = 56.834029
$ws.Range("I5").Value = 0.03347881112463321
$ws.Range("J5").Value = 0.03347881112463321
$ws.Range("M5").Value = 2.261165333333333
$ws.Range("N5").Value = 6.783496
$ws.Range("O5").Value = 0.04197723558351375
$ws.Range("P5").Value = 0.04197723558351374
$ws.Range("Q5").Value = 42.83704537615377
$ws.Range("R5").Value = 385.533408385384
$ws.Range("S5").Value = 0.001405347941634689
$ws.Range("T5").Value = 0.001405347941634689
$ws.Range("I6").Value = 0.4812547190371557
$ws.Range("J6").Value = 0.4812547190371557
$ws.Range("M6").Value = 15.03663066666667
$ws.Range("N6").Value = 45.109892
$ws.Range("O6").Value = 0.279146411176606
$ws.Range("P6").Value = 0.279146411176606
$ws.Range("Q6").Value = 4094.894303139445
$ws.Range("R6").Value = 36854.048728255
$ws.Range("S6").Value = 0.1343405276810279
$ws.Range("T6").Value = 0.1343405276810278
$ws.Range("I7").Value = 0.4812547190371557
$ws.Range("J7").Value = 0.4812547190371557
$ws.Range("O7").Value = 0.2673306493381863
$ws.Range("P7").Value = 0.2673306493381863
$ws.Range("S7").Value = 0.1286541365372692
$ws.Range("T7").Value = 0.1286541365372692
$ws.Range("I8").Value = 0.4812547190371557
$ws.Range("J8").Value = 0.4812547190371557
$ws.Range("M8").Value = 22.16851266666667
$ws.Range("N8").Value = 66.505538
$ws.Range("O8").Value = 0.411545703901694
$ws.Range("P8").Value = 0.411545703901694
$ws.Range("Q8").Value = 6037.104870111945
$ws.Range("R8").Value = 54333.9438310075
$ws.Range("S8").Value = 0.1980583121021582
$ws.Range("T8").Value = 0.1980583121021582
$ws.Range("I9").Value = 0.4812547190371557
$ws.Range("J9").Value = 0.4812547190371557
$ws.Range("M9").Value = 2.261165333333333
$ws.Range("N9").Value = 6.783496
$ws.Range("O9").Value = 0.04197723558351375
$ws.Range("P9").Value = 0.04197723558351374
$ws.Range("Q9").Value = 615.7784444655555
$ws.Range("R9").Value = 5542.006000189999
$ws.Range("S9").Value = 0.0202017427167004
$ws.Range("T9").Value = 0.0202017427167004
$ws.Range("G10").Value = 271.928284
$ws.Range("H10").Value = 815.784852
$ws.Range("I10").Value = 0.4805484928482698
$ws.Range("J10").Value = 0.4805484928482698
$ws.Range("M10").Value = 15.03663066666667
$ws.Range("N10").Value = 45.109892
$ws.Range("O10").Value = 0.279146411176606
$ws.Range("P10").Value = 0.279146411176606
$ws.Range("Q10").Value = 4088.885174328443
$ws.Range("R10").Value = 36799.96656895598
$ws.Range("S10").Value = 0.1341433871749215
$ws.Range("T10").Value = 0.1341433871749214
$ws.Range("G11").Value = 271.928284
$ws.Range("H11").Value = 815.784852
$ws.Range("I11").Value = 0.4805484928482698
$ws.Range("J11").Value = 0.4805484928482698
$ws.Range("O11").Value = 0.2673306493381863
$ws.Range("P11").Value = 0.2673306493381863
$ws.Range("Q11").Value = 3915.81007298335
$ws.Range("R11").Value = 35242.29065685015
$ws.Range("S11").Value = 0.1284653406316147
$ws.Range("T11").Value = 0.1284653406316147
$ws.Range("G12").Value = 271.928284
$ws.Range("H12").Value = 815.784852
$ws.Range("I12").Value = 0.4805484928482698
$ws.Range("J12").Value = 0.4805484928482698
$ws.Range("M12").Value = 22.16851266666667
$ws.Range("N12").Value = 66.505538
$ws.Range("O12").Value = 0.411545703901694
$ws.Range("P12").Value = 0.411545703901694
$ws.Range("Q12").Value = 6028.245608278931
$ws.Range("R12").Value = 54254.21047451038
$ws.Range("S12").Value = 0.1977676677481394
$ws.Range("T12").Value = 0.1977676677481393
$ws.Range("G13").Value = 271.928284
$ws.Range("H13").Value = 815.784852
$ws.Range("I13").Value = 0.4805484928482698
$ws.Range("J13").Value = 0.4805484928482698
$ws.Range("M13").Value = 2.261165333333333
$ws.Range("N13").Value = 6.783496
$ws.Range("O13").Value = 0.04197723558351375
$ws.Range("P13").Value = 0.04197723558351374
$ws.Range("Q13").Value = 614.8748089336212
$ws.Range("R13").Value = 5533.873280402591
$ws.Range("S13").Value = 0.02017209729359429
$ws.Range("T13").Value = 0.02017209729359429
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 2.669764666666667
$ws.Range("H14").Value = 8.009294000000001
$ws.Range("I14").Value = 0.004717976989941326
$ws.Range("J14").Value = 0.004717976989941326
$ws.Range("M14").Value = 15.03663066666667
$ws.Range("N14").Value = 45.109892
$ws.Range("O14").Value = 0.279146411176606
$ws.Range("P14").Value = 0.279146411176606
$ws.Range("Q14").Value = 40.14426525958311
$ws.Range("R14").Value = 361.298387336248
$ws.Range("S14").Value = 0.001317006344755928
$ws.Range("T14").Value = 0.001317006344755927
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2.669764666666667
$ws.Range("H15").Value = 8.009294000000001
$ws.Range("I15").Value = 0.004717976989941326
$ws.Range("J15").Value = 0.004717976989941326
$ws.Range("O15").Value = 0.2673306493381863
$ws.Range("P15").Value = 0.2673306493381863
$ws.Range("Q15").Value = 38.44503124297423
$ws.Range("R15").Value = 346.0052811867681
$ws.Range("S15").Value = 0.001261259852283636
$ws.Range("T15").Value = 0.001261259852283636
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 2.669764666666667
$ws.Range("H16").Value = 8.009294000000001
$ws.Range("I16").Value = 0.004717976989941326
$ws.Range("J16").Value = 0.004717976989941326
$ws.Range("M16").Value = 22.16851266666667
$ws.Range("N16").Value = 66.505538
$ws.Range("O16").Value = 0.411545703901694
$ws.Range("P16").Value = 0.411545703901694
$ws.Range("Q16").Value = 59.18471183001913
$ws.Range("R16").Value = 532.6624064701721
$ws.Range("S16").Value = 0.001941663161317399
$ws.Range("T16").Value = 0.001941663161317399
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 2.669764666666667
$ws.Range("H17").Value = 8.009294000000001
$ws.Range("I17").Value = 0.004717976989941326
$ws.Range("J17").Value = 0.004717976989941326
$ws.Range("M17").Value = 2.261165333333333
$ws.Range("N17").Value = 6.783496
$ws.Range("O17").Value = 0.04197723558351375
$ws.Range("P17").Value = 0.04197723558351374
$ws.Range("Q17").Value = 6.036779312424889
$ws.Range("R17").Value = 54.331013811824
$ws.Range("S17").Value = 0.0001980476315843641
$ws.Range("T17").Value = 0.0001980476315843641
